$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.581.29"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "3.652.97"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "196.33"
$ws.Range("E5").Value = "  +9.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "581.43"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "3.649.29"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  +2.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.684"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +7.96%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.66"
$ws.Range("E12").Value = "  +9.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000297"
$ws.Range("E13").Value = "  +19.11%  "
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("D15").Value = "4.234.24"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "3.655.67"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  +4.76%  "
$ws.Range("D19").Value = "68.499.93"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.71"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.03"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.89"
$ws.Range("E23").Value = "  +27.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.41"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("E26").Value = "  +5.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.71"
$ws.Range("E27").Value = "  +5.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.88"
$ws.Range("E28").Value = "  +8.50%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.22"
$ws.Range("E30").Value = "  +23.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.24"
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.88"
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "691.45"
$ws.Range("E33").Value = "  +17.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.30"
$ws.Range("E34").Value = "  +4.55%  "
$ws.Range("E35").Value = "  +7.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.04"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.90"
$ws.Range("E37").Value = "  +4.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.427"
$ws.Range("E38").Value = "  +15.43%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "0.0₃0797"
$ws.Range("E40").Value = "  +9.96%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.239.45"
$ws.Range("E41").Value = "  +21.15%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  +4.61%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("E43").Value = "  +20.80%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.14"
$ws.Range("E44").Value = "  +14.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").Value = "  +32.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0426"
$ws.Range("E47").Value = "  +4.85%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("E49").Value = "  +9.61%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.12"
$ws.Range("E50").Value = "  +3.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.62"
$ws.Range("E51").Value = "  +6.22%  "
